$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value2 = 210.77
$ws.Range("I15").Value2 = 210.77
$ws.Range("K15").Value2 = 632.3100000000001
$ws.Range("M15").Value2 = -463.3100000000001

$ws.Range("H62").Value2 = 348901.28
$ws.Range("I62").Value2 = 503475.4
$ws.Range("J62").Value2 = 5403.222
$ws.Range("K62").Value2 = 503475.4
$ws.Range("L62").Value2 = 5403.222
$ws.Range("M62").Value2 = -502851.4
$ws.Range("N62").Value2 = -6651.222

$ws.Range("H65").Value2 = 348901.28
$ws.Range("I65").Value2 = 503475.4
$ws.Range("J65").Value2 = 5403.222
$ws.Range("K65").Value2 = 2517377
$ws.Range("L65").Value2 = 27016.11
$ws.Range("M65").Value2 = -2514257
$ws.Range("N65").Value2 = -33256.11

$ws.Range("H98").Value2 = 1934.0952
$ws.Range("I98").Value2 = 2085.7896
$ws.Range("J98").Value2 = 493
$ws.Range("K98").Value2 = 2085.7896
$ws.Range("L98").Value2 = 493
$ws.Range("M98").Value2 = -587.7896000000001
$ws.Range("N98").Value2 = -3489

$ws.Range("H116").Value2 = 4780.3125
$ws.Range("I116").Value2 = 4623.75
$ws.Range("J116").Value2 = 5250
$ws.Range("K116").Value2 = 4623.75
$ws.Range("L116").Value2 = 5250
$ws.Range("M116").Value2 = -1181.75
$ws.Range("N116").Value2 = -12134

$ws.Range("H121").Value2 = 515.75
$ws.Range("I121").Value2 = 565
$ws.Range("J121").Value2 = 509.84
$ws.Range("K121").Value2 = 1695
$ws.Range("L121").Value2 = 1529.52
$ws.Range("M121").Value2 = 52
$ws.Range("N121").Value2 = -5023.52

$ws.Range("H122").Value2 = 1934.0952
$ws.Range("I122").Value2 = 2085.7896
$ws.Range("J122").Value2 = 493
$ws.Range("K122").Value2 = 6257.3688
$ws.Range("L122").Value2 = 1479
$ws.Range("M122").Value2 = -3807.3688
$ws.Range("N122").Value2 = -6379

$ws.Range("H123").Value2 = 46717.777
$ws.Range("J123").Value2 = 46717.777
$ws.Range("L123").Value2 = 46717.777
$ws.Range("N123").Value2 = -56517.777

$ws.Range("H141").Value2 = 2720
$ws.Range("I141").Value2 = 3933.3333
$ws.Range("K141").Value2 = 11799.9999
$ws.Range("M141").Value2 = -6619.999899999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value2 = 5173.5
$ws.Range("I32").Value2 = 4972.959
$ws.Range("J32").Value2 = 15000
$ws.Range("K32").Value2 = 4972.959
$ws.Range("L32").Value2 = 15000
$ws.Range("M32").Value2 = -4685.959
$ws.Range("N32").Value2 = -15574

$ws.Range("H74").Value2 = 1072.875
$ws.Range("I74").Value2 = 1084.2759
$ws.Range("J74").Value2 = 962.6667
$ws.Range("K74").Value2 = 1084.2759
$ws.Range("L74").Value2 = 962.6667
$ws.Range("M74").Value2 = -210.2759000000001
$ws.Range("N74").Value2 = -2710.6667

$ws.Range("H77").Value2 = 1072.875
$ws.Range("I77").Value2 = 1084.2759
$ws.Range("J77").Value2 = 962.6667
$ws.Range("K77").Value2 = 5421.379500000001
$ws.Range("L77").Value2 = 4813.3335
$ws.Range("M77").Value2 = -1053.379500000001
$ws.Range("N77").Value2 = -13549.3335

$ws.Range("H102").Value2 = 7001.6665
$ws.Range("I102").Value2 = 4910.909
$ws.Range("J102").Value2 = 30000
$ws.Range("K102").Value2 = 4910.909
$ws.Range("L102").Value2 = 30000
$ws.Range("M102").Value2 = -3288.909
$ws.Range("N102").Value2 = -33244

$ws.Range("H122").Value2 = 1244.7142
$ws.Range("I122").Value2 = 1052.6
$ws.Range("J122").Value2 = 1725
$ws.Range("K122").Value2 = 3157.8
$ws.Range("L122").Value2 = 5175
$ws.Range("M122").Value2 = -707.7999999999997
$ws.Range("N122").Value2 = -10075

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value2 = 45457644
$ws.Range("I86").Value2 = 71431250
$ws.Range("J86").Value2 = 3849
$ws.Range("K86").Value2 = 71431250
$ws.Range("L86").Value2 = 3849
$ws.Range("M86").Value2 = -71430127
$ws.Range("N86").Value2 = -6095

$ws.Range("H89").Value2 = 45457644
$ws.Range("I89").Value2 = 71431250
$ws.Range("J89").Value2 = 3849
$ws.Range("K89").Value2 = 357156250
$ws.Range("L89").Value2 = 19245
$ws.Range("M89").Value2 = -357150634
$ws.Range("N89").Value2 = -30477

$ws.Range("H92").Value2 = 20000
$ws.Range("J92").Value2 = 20000
$ws.Range("L92").Value2 = 20000
$ws.Range("N92").Value2 = -24992

$ws.Range("H107").Value2 = 4505.853
$ws.Range("I107").Value2 = 606.21875
$ws.Range("K107").Value2 = 606.21875
$ws.Range("M107").Value2 = 1313.78125

$ws.Range("H125").Value2 = 0
$ws.Range("J125").Value2 = 0
$ws.Range("L125").Value2 = 0
$ws.Range("N125").ClearContents()

$ws.Range("H126").Value2 = 57440
$ws.Range("J126").Value2 = 57440
$ws.Range("L126").Value2 = 57440
$ws.Range("N126").Value2 = -67320

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value2 = 38847.32
$ws.Range("I31").Value2 = 3201.087
$ws.Range("J31").Value2 = 202820
$ws.Range("K31").Value2 = 3201.087
$ws.Range("L31").Value2 = 202820
$ws.Range("M31").Value2 = -2906.087
$ws.Range("N31").Value2 = -203410

$ws.Range("H34").Value2 = 38847.32
$ws.Range("I34").Value2 = 3201.087
$ws.Range("J34").Value2 = 202820
$ws.Range("K34").Value2 = 3201.087
$ws.Range("L34").Value2 = 202820
$ws.Range("M34").Value2 = -2999.087
$ws.Range("N34").Value2 = -203224

$ws.Range("H58").Value2 = 2376.0725
$ws.Range("I58").Value2 = 861.9524
$ws.Range("K58").Value2 = 861.9524
$ws.Range("M58").Value2 = -658.9524

$ws.Range("H136").Value2 = 2376.0725
$ws.Range("I136").Value2 = 861.9524
$ws.Range("K136").Value2 = 2585.8572
$ws.Range("M136").Value2 = -35.85719999999992

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H97").Value2 = 369.23077
$ws.Range("I97").Value2 = 359.0909
$ws.Range("J97").Value2 = 425
$ws.Range("K97").Value2 = 1077.2727
$ws.Range("L97").Value2 = 1275
$ws.Range("M97").Value2 = -581.2727
$ws.Range("N97").Value2 = -2267

$ws.Range("H132").Value2 = 870.875
$ws.Range("I132").Value2 = 852.4286
$ws.Range("K132").Value2 = 7671.8574
$ws.Range("M132").Value2 = -5141.8574

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value2 = 4309.773
$ws.Range("I70").Value2 = 3855.3333
$ws.Range("J70").Value2 = 4624.385
$ws.Range("K70").Value2 = 3855.3333
$ws.Range("L70").Value2 = 4624.385
$ws.Range("M70").Value2 = -3585.3333
$ws.Range("N70").Value2 = -5164.385

$ws.Range("H73").Value2 = 4309.773
$ws.Range("I73").Value2 = 3855.3333
$ws.Range("J73").Value2 = 4624.385
$ws.Range("K73").Value2 = 3855.3333
$ws.Range("L73").Value2 = 4624.385
$ws.Range("M73").Value2 = -2919.3333
$ws.Range("N73").Value2 = -6496.385

$ws.Range("H97").Value2 = 3469.2307
$ws.Range("I97").Value2 = 3130
$ws.Range("J97").Value2 = 4600
$ws.Range("K97").Value2 = 3130
$ws.Range("L97").Value2 = 4600
$ws.Range("M97").Value2 = -2634
$ws.Range("N97").Value2 = -5592

$ws.Range("H132").Value2 = 2149.8293
$ws.Range("I132").Value2 = 1875.5714
$ws.Range("J132").Value2 = 3749.6667
$ws.Range("K132").Value2 = 5626.7142
$ws.Range("L132").Value2 = 11249.0001
$ws.Range("M132").Value2 = -3096.7142
$ws.Range("N132").Value2 = -16309.0001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value2 = 16033646
$ws.Range("I100").Value2 = 18705420
$ws.Range("J100").Value2 = 3000
$ws.Range("K100").Value2 = 18705420
$ws.Range("L100").Value2 = 3000
$ws.Range("M100").Value2 = -18704879
$ws.Range("N100").Value2 = -4082

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H105").Value2 = 31500
$ws.Range("J105").Value2 = 31500
$ws.Range("L105").Value2 = 31500
$ws.Range("N105").Value2 = -38488

$ws.Range("H107").Value2 = 565.4545000000001
$ws.Range("J107").Value2 = 425.375
$ws.Range("L107").Value2 = 1276.125
$ws.Range("N107").Value2 = -5116.125

$ws.Range("H122").Value2 = 24391208
$ws.Range("I122").Value2 = 27027968
$ws.Range("K122").Value2 = 81083904
$ws.Range("M122").Value2 = -81081454

$ws.Range("H132").Value2 = 1564.6923
$ws.Range("I132").Value2 = 1310.5
$ws.Range("J132").Value2 = 2136.625
$ws.Range("K132").Value2 = 3931.5
$ws.Range("L132").Value2 = 6409.875
$ws.Range("M132").Value2 = -1401.5
$ws.Range("N132").Value2 = -11469.875
